$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (row 53) - data shifted, one fewer forecast point
$ws.Rows("53").Delete()

# Update date (col A), year columns (B,D) and forecast ratio (C,E) values
$ws.Range("A2").Value = 39583
$ws.Range("B2").Value = 2008
$ws.Range("D2").Value = 2009
$ws.Range("E2").Value = 1.003756253906252
$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 0.8212989654785341
$ws.Range("A4").Value = 39948
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = 1.287693099940079
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 1.224010362214401
$ws.Range("A5").Value = 40130
$ws.Range("B5").Value = 2009
$ws.Range("C5").Value = 1.218009596270675
$ws.Range("D5").Value = 2010
$ws.Range("E5").Value = 1.183007486132071
$ws.Range("A6").Value = 40310
$ws.Range("B6").Value = 2010
$ws.Range("C6").Value = 0.5167526861706184
$ws.Range("D6").Value = 2011
$ws.Range("E6").Value = 0.9718821796794952
$ws.Range("A7").Value = 40494
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = 0.5544720893820188
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = 1.015842920196763
$ws.Range("A8").Value = 40676
$ws.Range("B8").Value = 2011
$ws.Range("C8").Value = 1.187829657075357
$ws.Range("D8").Value = 2012
$ws.Range("E8").Value = 1.00065194548169
$ws.Range("A9").Value = 40862
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 1.173294700162031
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 0.9092565586104273
$ws.Range("A10").Value = 41044
$ws.Range("B10").Value = 2012
$ws.Range("C10").Value = 1.071158385438342
$ws.Range("D10").Value = 2013
$ws.Range("E10").Value = 0.912403143334517
$ws.Range("A11").Value = 41228
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = 1.180518841971723
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 1.236730309040235
$ws.Range("A12").Value = 41409
$ws.Range("B12").Value = 2013
$ws.Range("C12").Value = 0.9553801317191413
$ws.Range("D12").Value = 2014
$ws.Range("E12").Value = 1.066801818459595
$ws.Range("A13").Value = 41592
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = 0.9276272455014611
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 1.029015928490629
$ws.Range("A14").Value = 41774
$ws.Range("B14").Value = 2014
$ws.Range("C14").Value = 1.190496724073231
$ws.Range("D14").Value = 2015
$ws.Range("E14").Value = 1.154811676806311
$ws.Range("A15").Value = 41957
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = 1.265990289415564
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 1.358148715145191
$ws.Range("A16").Value = 42137
$ws.Range("B16").Value = 2015
$ws.Range("C16").Value = 1.5464392869869
$ws.Range("D16").Value = 2016
$ws.Range("E16").Value = 1.247870081683522
$ws.Range("A17").Value = 42321
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 1.642047742738506
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 1.528208222695326
$ws.Range("A18").Value = 42503
$ws.Range("B18").Value = 2016
$ws.Range("C18").Value = 1.701952652941463
$ws.Range("D18").Value = 2017
$ws.Range("E18").Value = 1.637918813512695
$ws.Range("A19").Value = 42689
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = 1.66194179127146
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 1.634928000057778
$ws.Range("A20").Value = 42867
$ws.Range("B20").Value = 2017
$ws.Range("C20").Value = 1.580693894992691
$ws.Range("D20").Value = 2018
$ws.Range("E20").Value = 1.610567777412109
$ws.Range("A21").Value = 43053
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 1.609733807897773
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 1.67176973076042
$ws.Range("A22").Value = 43145
$ws.Range("B22").Value = 2018
$ws.Range("C22").Value = 1.651937828695615
$ws.Range("D22").Value = 2019
$ws.Range("E22").Value = 1.63821551487775
$ws.Range("A23").Value = 43235
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 1.646565058924154
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 1.636439239090515
$ws.Range("A24").Value = 43326
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 1.642460763882414
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = 1.62621273827539
$ws.Range("A25").Value = 43418
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 1.641178243814534
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 1.603287858019664
$ws.Range("A26").Value = 43510
$ws.Range("B26").Value = 2019
$ws.Range("C26").Value = 1.504616869537312
$ws.Range("D26").Value = 2020
$ws.Range("E26").Value = 1.599505522959732
$ws.Range("A27").Value = 43600
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 1.619750436871126
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 1.669486277487398
$ws.Range("A28").Value = 43691
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = 1.242963308065193
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 1.082447181878954
$ws.Range("A29").Value = 43783
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = 1.183163144818633
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = 0.8408455317168162
$ws.Range("A30").Value = 43875
$ws.Range("B30").Value = 2020
$ws.Range("C30").Value = 0.6216637650511503
$ws.Range("D30").Value = 2021
$ws.Range("E30").Value = 1.126729649114599
$ws.Range("A31").Value = 43966
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = 0.2954364073068261
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = 0.8326407735962826
$ws.Range("A32").Value = 44068
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = -3.662861831460751
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = -1.110565553434917
$ws.Range("A33").Value = 44159
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = -3.662861831460751
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = -1.875058665585216
$ws.Range("A34").Value = 44251
$ws.Range("B34").Value = 2021
$ws.Range("C34").Value = -3.604628722764358
$ws.Range("D34").Value = 2022
$ws.Range("E34").Value = -2.856219939917704
$ws.Range("A35").Value = 44341
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = -1.564297238929013
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = 0.1548119563699935
$ws.Range("A36").Value = 44432
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = 0.1010915562932313
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = 6.64637963114707
$ws.Range("A37").Value = 44525
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = 0.1010915562932313
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = 5.03478667886097
$ws.Range("A38").Value = 44617
$ws.Range("B38").Value = 2022
$ws.Range("C38").Value = 5.220550987750228
$ws.Range("D38").Value = 2023
$ws.Range("E38").Value = 1.043506288584606
$ws.Range("A39").Value = 44706
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = 5.937304773291885
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = 2.150399152794202
$ws.Range("A40").Value = 44798
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 5.793673192389748
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 1.728278600643907
$ws.Range("A41").Value = 44890
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = 5.793673192389748
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 2.399708479013141
$ws.Range("A42").Value = 44981
$ws.Range("B42").Value = 2023
$ws.Range("C42").Value = 0.156542203858212
$ws.Range("D42").Value = 2024
$ws.Range("E42").Value = 2.138492443986739
$ws.Range("A43").Value = 45071
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = -0.2621830498131694
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 1.878976297039481
$ws.Range("A44").Value = 45163
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = -0.3788601787194756
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 2.181307424743695
$ws.Range("A45").Value = 45254
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = -0.3788601787194756
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = 0.8520283695166997
$ws.Range("A46").Value = 45345
$ws.Range("B46").Value = 2024
$ws.Range("C46").Value = 0.005756553697899847
$ws.Range("D46").Value = 2025
$ws.Range("E46").Value = -0.1096192596443557
$ws.Range("A47").Value = 45436
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = 0.0512320434504332
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = 0.2660756331863467
$ws.Range("A48").Value = 45534
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = 0.05771202657300911
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 0.2051249733294291
$ws.Range("A49").Value = 45618
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = 0.05771202657300911
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = 0.299857156820571
$ws.Range("A50").Value = 45713
$ws.Range("B50").Value = 2025
$ws.Range("C50").Value = 0.4540776569412763
$ws.Range("D50").Value = 2026
$ws.Range("E50").Value = -0.1815195499670796
$ws.Range("A51").Value = 45800
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = 0.5998844096825495
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = 0.1733734969819434
$ws.Range("A52").Value = 45891
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = 0.6062046309774693
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = 0.3877310837361314
